# The commit deletes data row 43 ("# pièce" = 41846) from the "données"
# sheet: all rows below it shift up by one, the last (now-empty) row 75
# is cleared, the used range shrinks by one row, and the AutoFilter /
# _FilterDatabase range shrink from row 75 to row 74 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 43 and shift everything below it up by one.
$ws.Rows.Item(43).Delete()

# Deleting the row does not automatically resize the worksheet AutoFilter,
# so drop it and reapply over the new (one-row-shorter) range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:L74").AutoFilter()

# Keep the workbook-level _xlnm._FilterDatabase defined name in sync with
# the new AutoFilter range.
try {
    $n = $wb.Names.Item("données!_FilterDatabase")
} catch {
    $n = $wb.Names.Item(1)
}
$n.RefersTo = "=données!`$A`$1:`$L`$74"

# Restore the view: the user's selection ends up on A43 after the delete,
# scrolled so row 60 is at the top of the window.
$excel.ActiveWindow.ScrollRow = 60
[void]$ws.Range("A43").Select()
